$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-02 Thursday" "2025-10-03 Friday"

Replace-Text "793×2=" "859×7="
Replace-Text "190×7=" "528×5="
Replace-Text "692×7=" "551×5="
Replace-Text "450×3=" "458×2="
Replace-Text "905×2=" "818×5="

Replace-Text "595×2=" "625×3="
Replace-Text "473×7=" "476×3="
Replace-Text "442×5=" "345×3="
Replace-Text "338×3=" "827×8="
Replace-Text "471×2=" "872×3="

Replace-Text "767×7=" "812×3="
Replace-Text "672×2=" "365×8="
Replace-Text "186×9=" "885×7="
Replace-Text "663×3=" "481×8="
Replace-Text "266×8=" "683×7="

Replace-Text "591×8=" "102×5="
Replace-Text "862×3=" "659×2="
Replace-Text "378×6=" "988×9="
Replace-Text "921×6=" "624×3="
Replace-Text "290×6=" "665×3="

Replace-Text "512×7=" "154×9="
Replace-Text "138×7=" "838×4="
Replace-Text "850×6=" "365×9="
Replace-Text "313×2=" "465×2="
Replace-Text "986×2=" "470×7="
